$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing existing rows 47:56 down to 48:57
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with data
$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(47, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value = 45211
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(47, 6).Value = 100112044
$ws.Cells.Item(47, 7).Value = "Perejil"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 350
$ws.Cells.Item(47, 11).Value = 800
$ws.Cells.Item(47, 12).Value = 1000
$ws.Cells.Item(47, 13).Value = 914
$ws.Cells.Item(47, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 609
$ws.Cells.Item(47, 17).Value = 1.5
$ws.Cells.Item(47, 18).Value = "Hortaliza"
